# Swap the two theme color schemes used by this deck ("Office Theme" <->
# "Integral" / "Red Violet") so the slide master's theme now carries the
# colours that used to belong to the unused Office Theme part, matching the
# target commit (which physically swapped the contents of
# ppt/theme/theme1.xml and ppt/theme/theme2.xml).
#
# Only the colour scheme differs between the two theme parts in this deck
# (font scheme and format scheme are byte-identical), so re-pointing the
# slide master's 12 theme colours to the "Office Theme" palette reproduces
# the visible effect of that swap through the PowerPoint object model.

function Get-RGBValue($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Theme color order exposed by ThemeColorScheme is:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

$tcs.Item(1).RGB  = Get-RGBValue 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = Get-RGBValue 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = Get-RGBValue 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = Get-RGBValue 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = Get-RGBValue 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = Get-RGBValue 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = Get-RGBValue 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = Get-RGBValue 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = Get-RGBValue 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = Get-RGBValue 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = Get-RGBValue 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = Get-RGBValue 0x95 0x4F 0x72   # folHlink 954F72
